$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.380.08"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.610.12"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "2.608.74"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "3.090.83"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  -4.70%  "
$ws.Range("D17").Value = "67.296.39"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "2.612.14"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "367.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "66.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "579.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "0.0₃0988"
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -6.31%  "
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "0.0₆0287"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.26%  "
